# Workbook "diccionario.xlsx" edit:
#  - add a new worksheet "tbatraspasos" describing the xml reception/transfer
#    module (commit: "se agrego el modulo de cargar xml de recepcion de
#    transferencia de sucursal")
#  - tweak the existing "tbakardez" sheet: extend the note in C2 with the new
#    "I=Traspaso Ingreso" movement code and grow the row to fit it.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1) Add the new sheet right after "tbakardez" -------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "tbatraspasos"

# Header row
$ws2.Range("A1").Value = "Campo"
$ws2.Range("B1").Value = "Tipo de campo"
$ws2.Range("C1").Value = "Notas"

# Data row
$ws2.Range("A2").Value = "tipo"
$ws2.Range("B2").Value = "varchar"
$ws2.Range("C2").Value = "E= Envio,  R= Recepcion"

# Column C is a bit wider on this sheet
$ws2.Columns.Item(3).ColumnWidth = 12.6

# Row 2 needs to be taller to show the wrapped note
$ws2.Rows.Item(2).RowHeight = 60

# --- 2) Formatting for B2: a dedicated Arial 10pt black font --------------
# Build the font on a scratch named style first (keeps font creation to a
# minimum) and apply it, then drop the named style again so the workbook
# ends up with only the "Normal" cell style, same as before the edit.
$wb.Styles.Add("MiEstiloTemp") | Out-Null
$tmpStyle = $wb.Styles.Item("MiEstiloTemp")
$tmpStyle.Font.Name  = "Arial"
$tmpStyle.Font.Size  = 10
$tmpStyle.Font.Color = 0
$ws2.Range("B2").Style = "MiEstiloTemp"
$wb.Styles.Item("MiEstiloTemp").Delete()

# --- 3) Formatting for C2: wrap text + vertical-top (same look as the ----
#        "Notas" column on the first sheet)
$ws2.Range("C2").VerticalAlignment = -4160   # xlTop
$ws2.Range("C2").WrapText = $true

# Leave the new sheet's selection on C3, like a freshly filled-in sheet
$ws2.Range("C3").Select() | Out-Null

# --- 4) Update "tbakardez": append the new movement code to the note -----
$ws1.Range("C2").Value = "E= Entrada, S= salida, A= Ajuste, V= Venta, T= Traspaso, C= Compra, I=Traspaso Ingreso"
$ws1.Rows.Item(2).RowHeight = 135
$ws1.Range("C2").Select() | Out-Null

# --- 5) Restore "tbakardez" as the active/selected sheet ------------------
$ws1.Activate()
